$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the polyAIsolationProtocol column (G) for all 32 data rows:
# "NEBNextPoly(A)E7490L" -> "E7490L"
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 7).Value = "E7490L"
}

# Update the selection to reflect the new active cell/range
$ws.Range("G3:G33").Select()

$wb.Save()
